$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '67.703.54'
Set-TextValue 'E2' '  +0.46%  '
Set-TextValue 'D3' '2.495.96'
Set-TextValue 'E3' '  -2.38%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '590.72'
Set-TextValue 'E5' '  -0.42%  '
Set-TextValue 'D6' '174.00'
Set-TextValue 'E6' '  +0.25%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'E8' '  -1.18%  '
Set-TextValue 'D9' '2.496.12'
Set-TextValue 'E9' '  -2.38%  '
Set-TextValue 'D10' '0.139'
Set-TextValue 'E10' '  -0.23%  '
Set-TextValue 'E11' '  +1.68%  '
Set-TextValue 'D12' '5.10'
Set-TextValue 'E12' '  -1.36%  '
Set-TextValue 'D13' '0.342'
Set-TextValue 'E13' '  -2.79%  '
Set-TextValue 'D14' '26.30'
Set-TextValue 'E14' '  -3.15%  '
Set-TextValue 'D15' '2.948.04'
Set-TextValue 'E15' '  -2.30%  '
Set-TextValue 'E16' '  -1.41%  '
Set-TextValue 'D17' '67.678.00'
Set-TextValue 'D18' '2.483.00'
Set-TextValue 'E18' '  -2.85%  '
Set-TextValue 'D19' '11.77'
Set-TextValue 'E19' '  +3.12%  '
Set-TextValue 'D20' '7.99'
Set-TextValue 'E20' '  -0.90%  '
Set-TextValue 'D21' '365.08'
Set-TextValue 'E21' '  +2.30%  '
Set-TextValue 'D22' '4.13'
Set-TextValue 'E22' '  -2.58%  '
Set-TextValue 'D23' '4.56'
Set-TextValue 'E23' '  -2.74%  '
Set-TextValue 'D24' '71.38'
Set-TextValue 'E24' '  +1.65%  '
Set-TextValue 'E25' '  +0.09%  '
Set-TextValue 'E26' '  -5.68%  '
Set-TextValue 'D27' '9.96'
Set-TextValue 'E27' '  -2.64%  '
Set-TextValue 'D28' '0.998'
Set-TextValue 'E28' '  -0.14%  '
Set-TextValue 'D29' '2.622.35'
Set-TextValue 'E29' '  -2.22%  '
Set-TextValue 'D30' '0.0₃0962'
Set-TextValue 'E30' '  -3.79%  '
Set-TextValue 'B31' 'InternetComputer(DFINITY)'
Set-TextValue 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D31' '8.32'
Set-TextValue 'E31' '  +1.12%  '
Set-TextValue 'B32' 'Bittensor'
Set-TextValue 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D32' '531.56'
Set-TextValue 'E32' '  -1.47%  '
Set-TextValue 'E33' '  -5.66%  '
Set-TextValue 'E34' '  -0.14%  '
Set-TextValue 'E35' '  +0.03%  '
Set-TextValue 'D36' '0.127'
Set-TextValue 'E36' '  -4.10%  '
Set-TextValue 'D37' '157.97'
Set-TextValue 'E37' '  -0.42%  '
Set-TextValue 'E38' '  -4.02%  '
Set-TextValue 'D39' '18.71'
Set-TextValue 'E39' '  -0.54%  '
Set-TextValue 'E40' '  +0.95%  '
Set-TextValue 'B41' 'PolygonEcosystemToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D41' '0.349'
Set-TextValue 'E41' '  -2.72%  '
Set-TextValue 'B42' 'Stacks'
Set-TextValue 'C42' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D42' '1.78'
Set-TextValue 'E42' '  -1.84%  '
Set-TextValue 'D43' '5.10'
Set-TextValue 'E43' '  -1.94%  '
Set-TextValue 'B44' 'USDe'
Set-TextValue 'C44' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D44' '1.00'
Set-TextValue 'E44' '  -0.03%  '
Set-TextValue 'B45' 'dogwifhat'
Set-TextValue 'C45' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D45' '2.51'
Set-TextValue 'E45' '  -1.27%  '
Set-TextValue 'D46' '145.20'
Set-TextValue 'E46' '  -3.97%  '
Set-TextValue 'D47' '3.68'
Set-TextValue 'E47' '  -1.42%  '
Set-TextValue 'B48' 'ARBITRUM'
Set-TextValue 'C48' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D48' '0.547'
Set-TextValue 'E48' '  -3.46%  '
Set-TextValue 'B49' 'BabyDogeCoin'
Set-TextValue 'C49' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D49' '0.0₆0274'
Set-TextValue 'E49' '  -2.89%  '
Set-TextValue 'D50' '1.69'
Set-TextValue 'E50' '  -2.12%  '
